$d = $word.ActiveDocument

$pairs = @(
    @("242÷2=121, 0", "856÷8=107, 0"),
    @("442÷9=49, 1", "329÷9=36, 5"),
    @("820÷8=102, 4", "834÷8=104, 2"),
    @("955÷5=191, 0", "880÷4=220, 0"),
    @("494÷8=61, 6", "933÷7=133, 2"),
    @("230÷4=57, 2", "329÷7=47, 0"),
    @("887÷5=177, 2", "602÷3=200, 2"),
    @("738÷4=184, 2", "986÷3=328, 2"),
    @("497÷7=71, 0", "646÷8=80, 6"),
    @("825÷8=103, 1", "611÷8=76, 3"),
    @("937÷2=468, 1", "491÷7=70, 1"),
    @("606÷8=75, 6", "285÷3=95, 0"),
    @("395÷5=79, 0", "177÷3=59, 0"),
    @("136÷5=27, 1", "558÷4=139, 2"),
    @("945÷6=157, 3", "843÷9=93, 6"),
    @("717÷5=143, 2", "130÷2=65, 0"),
    @("790÷2=395, 0", "521÷4=130, 1"),
    @("114÷8=14, 2", "379÷3=126, 1"),
    @("249÷2=124, 1", "418÷4=104, 2"),
    @("370÷2=185, 0", "383÷2=191, 1"),
    @("290÷5=58, 0", "899÷7=128, 3"),
    @("621÷5=124, 1", "792÷8=99, 0"),
    @("544÷6=90, 4", "696÷7=99, 3"),
    @("803÷2=401, 1", "802÷7=114, 4"),
    @("114÷2=57, 0", "958÷6=159, 4")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
